# Daily attendance processing - 2025-10-11 06:25:23
# Reverses the order of the comma-separated "Recorded By" entries in
# column G for the affected rows on the "Session Analysis Results" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# Rows in column G whose "Recorded By" list order needs to be reversed.
$rows = @(2,3,4,5,6,7,10,11,12,13,14,15,29,30,32,33,34,37,38,39,40,41,42,56,57,58,59,60,61,64,65,66,67,68,69,84,85,86,87,88,89,90,93,95,110,111,112,113,114,115,116,119,121,136,137,138,139,140,141,142,145,147)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = 7
    $val = $cell.Value2
    if ($null -ne $val -and $val -ne "") {
        $parts = $val -split ", "
        $n = $parts.Count
        $reversed = @()
        for ($i = $n - 1; $i -ge 0; $i--) {
            $reversed += $parts[$i]
        }
        $cell.Value = [string]::Join(", ", $reversed)
    }
}
